$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1425304360311941
$ws.Range("D2").Value = 0.04416460176246062
$ws.Range("E2").Value = 0.4396191691992186
$ws.Range("F2").Value = 1.528290297222981
$ws.Range("G2").Value = 1.673767081463382
$ws.Range("H2").Value = 1.010542361786406
$ws.Range("N2").Value = 1.995241137946778

# Row 3
$ws.Range("B3").Value = 0.1329091933101978
$ws.Range("D3").Value = 0.03830832642090343
$ws.Range("E3").Value = 0.3823479279035951
$ws.Range("F3").Value = 1.363286054659909
$ws.Range("G3").Value = 1.472460773835678
$ws.Range("H3").Value = 0.9231199012414208
$ws.Range("N3").Value = 1.846526605564804

# Row 4
$ws.Range("B4").Value = 0.1270661492628875
$ws.Range("D4").Value = 0.03472637085036467
$ws.Range("E4").Value = 0.3474053453013823
$ws.Range("F4").Value = 1.262945067803372
$ws.Range("G4").Value = 1.349747030403535
$ws.Range("H4").Value = 0.8702008080870769
$ws.Range("N4").Value = 1.755520437593759

# Row 5
$ws.Range("B5").Value = 0.1247015078378411
$ws.Range("D5").Value = 0.03326964694228707
$ws.Range("E5").Value = 0.333216280922386
$ws.Range("F5").Value = 1.222288067622259
$ws.Range("G5").Value = 1.299949458491369
$ws.Range("H5").Value = 0.8488200623750686
$ws.Range("N5").Value = 1.71851688491077

# Row 6
$ws.Range("B6").Value = 0.1243098614095288
$ws.Range("D6").Value = 0.03302792226286044
$ws.Range("E6").Value = 0.3308630769826664
$ws.Range("F6").Value = 1.215550724464379
$ws.Range("G6").Value = 1.291692836908965
$ws.Range("H6").Value = 0.8452807431618226
$ws.Range("N6").Value = 1.712377599108464

# Row 7
$ws.Range("B7").Value = 0.1270341919676241
$ws.Range("D7").Value = 0.03470671369345268
$ws.Range("E7").Value = 0.3472137901143384
$ws.Range("F7").Value = 1.262395826126067
$ws.Range("G7").Value = 1.349074614279345
$ws.Range("H7").Value = 0.8699117225502562
$ws.Range("N7").Value = 1.755021054646249

# Row 8
$ws.Range("B8").Value = 0.1391997998655086
$ws.Range("D8").Value = 0.04214216942487781
$ws.Range("E8").Value = 0.4198225367166089
$ws.Range("F8").Value = 1.47118785694326
$ws.Range("G8").Value = 1.604163365212003
$ws.Range("H8").Value = 0.98023762357343
$ws.Range("N8").Value = 1.943904251618477

# Row 9
$ws.Range("B9").Value = 0.1635584291338716
$ws.Range("D9").Value = 0.05685744970638495
$ws.Range("E9").Value = 0.5642327170640442
$ws.Range("F9").Value = 1.888902018709359
$ws.Range("G9").Value = 2.112127782647462
$ws.Range("H9").Value = 1.202918706961327
$ws.Range("N9").Value = 2.316504389646809

# Row 10
$ws.Range("B10").Value = 0.1817500012201947
$ws.Range("D10").Value = 0.06778656051733378
$ws.Range("E10").Value = 0.6719446234360191
$ws.Range("F10").Value = 2.201659492699889
$ws.Range("G10").Value = 2.491042104895939
$ws.Range("H10").Value = 1.370841485233257
$ws.Range("N10").Value = 2.591327722520418

# Row 11
$ws.Range("B11").Value = 0.190087785889915
$ws.Range("D11").Value = 0.07279217337988086
$ws.Range("E11").Value = 0.7213816670721798
$ws.Range("F11").Value = 2.345399064368451
$ws.Range("G11").Value = 2.664885521213535
$ws.Range("H11").Value = 1.44827547428082
$ws.Range("N11").Value = 2.716530669508643

# Row 12
$ws.Range("B12").Value = 0.1932538308941787
$ws.Range("D12").Value = 0.07469321267035411
$ws.Range("E12").Value = 0.7401723167001393
$ws.Range("F12").Value = 2.400055170825851
$ws.Range("G12").Value = 2.730945692543116
$ws.Range("H12").Value = 1.477756451588164
$ws.Range("N12").Value = 2.763963193734014

# Row 13
$ws.Range("B13").Value = 0.19257158332789
$ws.Range("D13").Value = 0.0742835337146488
$ws.Range("E13").Value = 0.7361221926845474
$ws.Range("F13").Value = 2.388273756458403
$ws.Range("G13").Value = 2.71670796254466
$ws.Range("H13").Value = 1.471400018306724
$ws.Range("N13").Value = 2.753746911386088

# Row 14
$ws.Range("B14").Value = 0.1903480852984956
$ws.Range("D14").Value = 0.0729484585700817
$ws.Range("E14").Value = 0.7229261437163217
$ws.Range("F14").Value = 2.349891060477688
$ws.Range("G14").Value = 2.670315630292464
$ws.Range("H14").Value = 1.450697669861711
$ws.Range("N14").Value = 2.720432585385765

# Row 15
$ws.Range("B15").Value = 0.1889872547797466
$ws.Range("D15").Value = 0.07213142522000737
$ws.Range("E15").Value = 0.7148524956228215
$ws.Range("F15").Value = 2.326410295331868
$ws.Range("G15").Value = 2.641929398215893
$ws.Range("H15").Value = 1.438037758433723
$ws.Range("N15").Value = 2.700029152110062

# Row 16
$ws.Range("B16").Value = 0.1812063398626123
$ws.Range("D16").Value = 0.06746016902452823
$ws.Range("E16").Value = 0.6687232001523142
$ws.Range("F16").Value = 2.192296469740057
$ws.Range("G16").Value = 2.479712166869206
$ws.Range("H16").Value = 1.365802710240757
$ws.Range("N16").Value = 2.583148644060884

# Row 17
$ws.Range("B17").Value = 0.1764487851811509
$ws.Range("D17").Value = 0.06460364233217319
$ws.Range("E17").Value = 0.6405414693056741
$ws.Range("F17").Value = 2.110407240501445
$ws.Range("G17").Value = 2.38058685249996
$ws.Range("H17").Value = 1.321762217021217
$ws.Range("N17").Value = 2.511489453097226

# Row 18
$ws.Range("B18").Value = 0.1737182494248515
$ws.Range("D18").Value = 0.06296379515380579
$ws.Range("E18").Value = 0.6243729184997164
$ws.Range("F18").Value = 2.063443689084153
$ws.Range("G18").Value = 2.323710192705676
$ws.Range("H18").Value = 1.2965291561203
$ws.Range("N18").Value = 2.470290675910178

# Row 19
$ws.Range("B19").Value = 0.1727947553762306
$ws.Range("D19").Value = 0.06240909297251562
$ws.Range("E19").Value = 0.6189053326469747
$ws.Range("F19").Value = 2.047565729167701
$ws.Range("G19").Value = 2.304475839262807
$ws.Range("H19").Value = 1.288002240429705
$ws.Range("N19").Value = 2.456344689235891

# Row 20
$ws.Range("B20").Value = 0.176954627914796
$ws.Range("D20").Value = 0.06490739369498044
$ws.Range("E20").Value = 0.6435371872360207
$ws.Range("F20").Value = 2.119110209985365
$ws.Range("G20").Value = 2.391124544869683
$ws.Range("H20").Value = 1.326440215915113
$ws.Range("N20").Value = 2.519115905984449

# Row 21
$ws.Range("B21").Value = 0.1910009463247206
$ws.Range("D21").Value = 0.07334044737054057
$ws.Range("E21").Value = 0.726800189790751
$ws.Range("F21").Value = 2.361158763589458
$ws.Range("G21").Value = 2.683935813114317
$ws.Range("H21").Value = 1.45677408627796
$ws.Range("N21").Value = 2.730217288371648

# Row 22
$ws.Range("B22").Value = 0.2002316580427816
$ws.Range("D22").Value = 0.07888443913425647
$ws.Range("E22").Value = 0.7816282444840539
$ws.Range("F22").Value = 2.520670277776958
$ws.Range("G22").Value = 2.876650817821144
$ws.Range("H22").Value = 1.542882076632281
$ws.Range("N22").Value = 2.8683033949373

# Row 23
$ws.Range("B23").Value = 0.1953005061487829
$ws.Range("D23").Value = 0.07592231310724173
$ws.Range("E23").Value = 0.752325552932561
$ws.Range("F23").Value = 2.435410421870188
$ws.Range("G23").Value = 2.773666093207737
$ws.Range("H23").Value = 1.496837023020134
$ws.Range("N23").Value = 2.794595190493112

# Row 24
$ws.Range("B24").Value = 0.1767259219250121
$ws.Range("D24").Value = 0.06477006025150445
$ws.Range("E24").Value = 0.6421827193271525
$ws.Range("F24").Value = 2.115175238321626
$ws.Range("G24").Value = 2.386360107228029
$ws.Range("H24").Value = 1.324325024093469
$ws.Range("N24").Value = 2.515667989176677

# Row 25
$ws.Range("B25").Value = 0.1569161479393557
$ws.Range("D25").Value = 0.05285854641996934
$ws.Range("E25").Value = 0.524910646920091
$ws.Range("F25").Value = 1.774930254878313
$ws.Range("G25").Value = 1.973780395381823
$ws.Range("H25").Value = 1.14195305260256
$ws.Range("N25").Value = 2.215498394811675
